$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - styled the same as the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Timestamp values for rows 2-23
$times = @(
    "2021-10-05 10:50:18.374731",
    "2021-10-05 10:50:18.374741",
    "2021-10-05 10:50:18.374744",
    "2021-10-05 10:50:18.374747",
    "2021-10-05 10:50:18.374750",
    "2021-10-05 10:50:18.374753",
    "2021-10-05 10:50:18.374755",
    "2021-10-05 10:50:18.374758",
    "2021-10-05 10:50:18.374761",
    "2021-10-05 10:50:18.374763",
    "2021-10-05 10:50:18.374766",
    "2021-10-05 10:50:18.374768",
    "2021-10-05 10:50:18.374771",
    "2021-10-05 10:50:18.374773",
    "2021-10-05 10:50:18.374776",
    "2021-10-05 10:50:18.374778",
    "2021-10-05 10:50:18.374781",
    "2021-10-05 10:50:18.374784",
    "2021-10-05 10:50:18.374787",
    "2021-10-05 10:50:18.374789",
    "2021-10-05 10:50:18.374791",
    "2021-10-05 10:50:18.374794"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
